$wb = $excel.ActiveWorkbook
Write-Host "excel type test"
try { $wb.RunCommand("pivot.refresh", '{}') } catch { Write-Host "err1: $_" }
try { $excel.RunCommand("pivot.refresh", '{}') } catch { Write-Host "err2: $_" }
